$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.205.62"
$ws.Range("E2").Value = "'  -0.50%  "
$ws.Range("D3").Value = "'1.828.54"
$ws.Range("E3").Value = "'  -0.69%  "
$ws.Range("D4").Value = "'0.9987"
$ws.Range("E4").Value = "'  +0.03%  "
$ws.Range("D5").Value = "'236.68"
$ws.Range("E5").Value = "'  -1.32%  "
$ws.Range("D6").Value = "'0.6071"
$ws.Range("E7").Value = "'  +0.05%  "
$ws.Range("D8").Value = "'0.07098"
$ws.Range("E8").Value = "'  -4.82%  "
$ws.Range("D9").Value = "'0.2814"
$ws.Range("E9").Value = "'  -2.84%  "
$ws.Range("D10").Value = "'23.84"
$ws.Range("E10").Value = "'  -4.84%  "
$ws.Range("D11").Value = "'0.07665"
$ws.Range("E11").Value = "'  -0.80%  "
$ws.Range("D12").Value = "'1.825.18"
$ws.Range("E12").Value = "'  -0.76%  "
$ws.Range("E13").Value = "'  -2.90%  "
$ws.Range("E14").Value = "'  -2.34%  "
$ws.Range("D15").Value = "'0.6369"
$ws.Range("E15").Value = "'  -5.92%  "
$ws.Range("D16").Value = "'2.071.66"
$ws.Range("E16").Value = "'  -0.81%  "
$ws.Range("D17").Value = "'79.30"
$ws.Range("E17").Value = "'  -3.09%  "
$ws.Range("D18").Value = "'5.913"
$ws.Range("E18").Value = "'  -5.19%  "
$ws.Range("D19").Value = "'29.177.85"
$ws.Range("E19").Value = "'  -0.60%  "
$ws.Range("D20").Value = "'228.24"
$ws.Range("E20").Value = "'  -0.33%  "
$ws.Range("D21").Value = "'11.81"
$ws.Range("E21").Value = "'  -4.14%  "
$ws.Range("D22").Value = "'1.000"
$ws.Range("D23").Value = "'7.031"
$ws.Range("E23").Value = "'  -4.66%  "
$ws.Range("D24").Value = "'1.000"
$ws.Range("E24").Value = "'  +0.07%  "
$ws.Range("D25").Value = "'154.54"
$ws.Range("E25").Value = "'  -2.21%  "
$ws.Range("D26").Value = "'8.077"
$ws.Range("E26").Value = "'  -5.33%  "
$ws.Range("D27").Value = "'0.1297"
$ws.Range("E27").Value = "'  -3.82%  "
$ws.Range("D28").Value = "'16.59"
$ws.Range("E28").Value = "'  -4.81%  "
$ws.Range("D29").Value = "'1.480"
$ws.Range("E29").Value = "'  +1.82%  "
$ws.Range("D30").Value = "'0.06500"
$ws.Range("E30").Value = "'  -5.83%  "
$ws.Range("D31").Value = "'1.458"
$ws.Range("E31").Value = "'  -2.12%  "
$ws.Range("D32").Value = "'3.832"
$ws.Range("E32").Value = "'  -5.52%  "
$ws.Range("D33").Value = "'3.814"
$ws.Range("E33").Value = "'  -6.20%  "
$ws.Range("D34").Value = "'1.129"
$ws.Range("E34").Value = "'  -0.95%  "
$ws.Range("D35").Value = "'1.750"
$ws.Range("E35").Value = "'  -4.33%  "
$ws.Range("D36").Value = "'0.6504"
$ws.Range("E36").Value = "'  -6.90%  "
$ws.Range("D37").Value = "'2.552"
$ws.Range("E37").Value = "'  -1.17%  "
$ws.Range("D38").Value = "'2.755"
$ws.Range("E38").Value = "'  -2.40%  "
$ws.Range("D39").Value = "'1.215.90"
$ws.Range("E39").Value = "'  -1.75%  "
$ws.Range("D40").Value = "'0.01753"
$ws.Range("E40").Value = "'  -5.05%  "
$ws.Range("D41").Value = "'6.500"
$ws.Range("E41").Value = "'  -4.37%  "
$ws.Range("D42").Value = "'0.9316"
$ws.Range("E42").Value = "'  -0.94%  "
$ws.Range("D43").Value = "'0.9995"
$ws.Range("D44").Value = "'101.02"
$ws.Range("E44").Value = "'  -0.09%  "
$ws.Range("D45").Value = "'1.982.51"
$ws.Range("E45").Value = "'  +0.25%  "
$ws.Range("D46").Value = "'63.20"
$ws.Range("E46").Value = "'  -3.37%  "
$ws.Range("D47").Value = "'0.00000000117"
$ws.Range("E47").Value = "'  -1.49%  "
$ws.Range("B48").Value = "'RenderToken"
$ws.Range("C48").Value = "'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").Value = "'1.612"
$ws.Range("E48").Value = "'  -5.80%  "
$ws.Range("B49").Value = "'EnergySwap"
$ws.Range("C49").Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'8.578"
$ws.Range("E49").Value = "'  -4.21%  "
$ws.Range("E50").Value = "'  -5.56%  "
$ws.Range("D51").Value = "'0.05528"
$ws.Range("E51").Value = "'  -2.61%  "
